# "add data in testdata file"
# Login sheet gains a new email row (navi3@gmail.com) inserted right after the
# existing two addresses, the old trailing navi4@gmail.com row is dropped, and
# the Register sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")
$wsRegister = $wb.Worksheets.Item("Register")

# ---------------------------------------------------------------------------
# Login sheet: reshuffle the email column.
#   old: row2=navi1  row3=navi2  row4=navi3  row5=navi4
#   new: row2=navi3  row3=navi1  row4=navi2            (navi4 row removed)
# ---------------------------------------------------------------------------

# Drop all existing hyperlinks up front so none are left dangling once row 5
# disappears (individual hyperlink deletion isn't reliable, so clear + rebuild).
$wsLogin.Hyperlinks.Delete()

# Remove the last data row (used to hold navi4@gmail.com).
$wsLogin.Rows.Item(5).Delete() | Out-Null

# Write the reordered + new email addresses.
$wsLogin.Range("A2").Value = "navi3@gmail.com"
$wsLogin.Range("A3").Value = "navi1@gmail.com"
$wsLogin.Range("A4").Value = "navi2@gmail.com"

# Recreate the mailto hyperlinks for the three remaining rows.
$wsLogin.Hyperlinks.Add($wsLogin.Range("A2"), "mailto:navi3@gmail.com") | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("A3"), "mailto:navi1@gmail.com") | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("A4"), "mailto:navi2@gmail.com") | Out-Null

# Keep the usual hyperlink look on those cells.
$wsLogin.Range("A2").Style = "Hyperlink"
$wsLogin.Range("A3").Style = "Hyperlink"
$wsLogin.Range("A4").Style = "Hyperlink"

$wsLogin.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Register sheet: tidy the column widths and become the active sheet/tab.
# ---------------------------------------------------------------------------
$wsRegister.Columns.Item(1).ColumnWidth = 18
$wsRegister.Columns.Item(2).ColumnWidth = 15.333333333333334
$wsRegister.Columns.Item(3).ColumnWidth = 32
$wsRegister.Columns.Item(4).ColumnWidth = 19.333333333333332

$wsRegister.Activate() | Out-Null
$wsRegister.Range("A1:E3").Select() | Out-Null
